# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6
$ws1.Range("F3").Value = 300
$ws1.Range("F5").Value = 2548
$ws1.Range("F6").Value = 1838
$ws1.Range("F9").Value = 893
$ws1.Range("F10").Value = 175

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6
$ws4.Range("F3").Value = 300
$ws4.Range("F5").Value = 2548
$ws4.Range("F6").Value = 1838
$ws4.Range("F10").Value = 893
$ws4.Range("F11").Value = 175
